# Données_groupe_15 - fix data values ("logic problems" per commit message)
# and refresh the last-used-folder / window-size bookkeeping attributes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("données15")

# --- Cosmetic workbook metadata (best effort; harmless if the host ignores them) ---
# Window size Excel remembers for next time the workbook is opened.
$win = $wb.Windows.Item(1)
$win.Width  = 25800
$win.Height = 13200

# --- Data corrections ---
$ws.Range("A14").Value = 7.180000000000001
$ws.Range("C14").Value = 87
$ws.Range("A18").Value = 31.59
$ws.Range("C18").Value = 86
$ws.Range("A19").Value = 19.470000000000002
$ws.Range("C19").Value = 93
$ws.Range("A20").Value = 11.21
$ws.Range("C20").Value = 84
$ws.Range("A21").Value = 3.9699999999999998
$ws.Range("C21").Value = 91
$ws.Range("A22").Value = 33.47
$ws.Range("C22").Value = 86
$ws.Range("A23").Value = 43.230000000000004
$ws.Range("C23").Value = 74
$ws.Range("A24").Value = 18.81
$ws.Range("C24").Value = 90
$ws.Range("A26").Value = 15.8
$ws.Range("C26").Value = 79
$ws.Range("A27").Value = 6.710000000000001
$ws.Range("C27").Value = 92
$ws.Range("A30").Value = 3.42
$ws.Range("C30").Value = 39
$ws.Range("A31").Value = 46.760000000000005
$ws.Range("C31").Value = 91
$ws.Range("A32").Value = 27.529999999999998
$ws.Range("C32").Value = 70
$ws.Range("A34").Value = 17.66
$ws.Range("C34").Value = 88
$ws.Range("A35").Value = 3.95
$ws.Range("C35").Value = 90
$ws.Range("A36").Value = 22.99
$ws.Range("C36").Value = 85
$ws.Range("A37").Value = 14.580000000000002
$ws.Range("C37").Value = 55
$ws.Range("A39").Value = 7.35
$ws.Range("C39").Value = 79
$ws.Range("A40").Value = 72.43
$ws.Range("C40").Value = 87
$ws.Range("A42").Value = 4.52
$ws.Range("C42").Value = 92
$ws.Range("A43").Value = 27.1
$ws.Range("C43").Value = 93
$ws.Range("A45").Value = 5.09
$ws.Range("C45").Value = 92
$ws.Range("A46").Value = 7.46
$ws.Range("C46").Value = 78
$ws.Range("A47").Value = 47.79
$ws.Range("C47").Value = 83
$ws.Range("A48").Value = 51.29
$ws.Range("C48").Value = 92
$ws.Range("A49").Value = 22.770000000000003
$ws.Range("C49").Value = 80
$ws.Range("A50").Value = 16.31
$ws.Range("C50").Value = 58
$ws.Range("A53").Value = 12.02
$ws.Range("C53").Value = 65
$ws.Range("A55").Value = 8.91
$ws.Range("C55").Value = 62
$ws.Range("A56").Value = 5.92
$ws.Range("C56").Value = 93
$ws.Range("A57").Value = 14.69
$ws.Range("C57").Value = 84
$ws.Range("A61").Value = 2.97
$ws.Range("C61").Value = 90
$ws.Range("A63").Value = 10.97
$ws.Range("C63").Value = 46
$ws.Range("A64").Value = 30.19
$ws.Range("C64").Value = 93
$ws.Range("A65").Value = 15.06
$ws.Range("C65").Value = 70
$ws.Range("A66").Value = 11.55
$ws.Range("C66").Value = 79
$ws.Range("A67").Value = 4.66
$ws.Range("C67").Value = 78
$ws.Range("A68").Value = 17.91
$ws.Range("C68").Value = 70
$ws.Range("A69").Value = 5.9499999999999993
$ws.Range("C69").Value = 78
$ws.Range("A70").Value = 32.84
$ws.Range("C70").Value = 59
$ws.Range("A71").Value = 5.62
$ws.Range("C71").Value = 58
$ws.Range("A72").Value = 16.32
$ws.Range("C72").Value = 65
$ws.Range("A73").Value = 47.54
$ws.Range("C73").Value = 79
$ws.Range("A74").Value = 17.91
$ws.Range("C74").Value = 92
$ws.Range("A75").Value = 12.23
$ws.Range("C75").Value = 76
$ws.Range("A76").Value = 20.22
$ws.Range("C76").Value = 81
$ws.Range("A77").Value = 11.799999999999999
$ws.Range("C77").Value = 75
$ws.Range("A78").Value = 48.94
$ws.Range("C78").Value = 64
$ws.Range("A80").Value = 5.21
$ws.Range("C80").Value = 48
